# Word COM-interop script for the "logico_descritivo.docx" edit.
#
# The original first paragraph reads:
#   "funcionario (id_funcionario, nome, sobrenome, login, senha, cpf,
#    telefone, #cod_setor);"
# (built from runs "funcionario (" + "id_f" + [bookmark _GoBack] +
#  "uncionario" + ", nome, ... #cod_setor);").
#
# It is replaced by four paragraphs:
#   1. "Alunos: Guilherme Francisco e Yasmin Moraes"
#   2. "Turma: 2K"
#   3. (empty paragraph holding the _GoBack bookmark)
#   4. "funcionario (id_funcionario, nome, sobrenome, login, senha, cpf,
#       telefone, #cod_setor);" -- same visible text as before, but now
#      "id_f"+"uncionario" is a single "id_funcionario" run, and
#      spell-check proofErr markers wrap "funcionario" and
#      "id_funcionario".
#
# All other paragraphs in the document are left untouched.

$d = $word.ActiveDocument

# Locate the paragraph that starts with "funcionario (" (the one the
# diff rewrites) rather than assuming it is always Paragraphs(1).
$searchRange = $d.Content
$found = $searchRange.Find.Execute("funcionario (", $true, $false, $false,
                                    $false, $false, $true, 1, $false,
                                    "", 0)

if ($found) {
    $targetPara = $searchRange.Paragraphs(1)
} else {
    $targetPara = $d.Paragraphs(1)
}

$target = $targetPara.Range

$xml = "<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:pPr><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/></w:rPr><w:t>Alunos: Guilherme Francisco e Yasmin Moraes</w:t></w:r></w:p><w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:pPr><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/></w:rPr><w:t>Turma: 2K</w:t></w:r></w:p><w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:pPr><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/></w:rPr></w:pPr><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p><w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:pPr><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/></w:rPr></w:pPr><w:proofErr w:type=`"spellStart`"/><w:r><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/></w:rPr><w:t>funcionario</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/></w:rPr><w:t xml:space=`"preserve`"> (</w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:u w:val=`"single`"/></w:rPr><w:t>id_funcionario</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/></w:rPr><w:t>, nome, sobrenome, login, senha, cpf, telefone, #cod_setor);</w:t></w:r></w:p>"

$target.InsertXML($xml)

Write-Output ("Paragraphs after edit: " + $d.Paragraphs.Count)
